$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-46, columns I and J
$data = @{
    2  = @(10,10)
    3  = @(4,5)
    4  = @(9,9)
    5  = @(5,6)
    6  = @(9,9)
    7  = @(5,6)
    8  = @(7,8)
    9  = @(6,7)
    10 = @(9,9)
    11 = @(7,8)
    12 = @(6,6)
    13 = @(7,8)
    14 = @(8,9)
    15 = @(4,6)
    16 = @(7,9)
    17 = @(5,6)
    18 = @(3,6)
    19 = @(7,8)
    20 = @(6,7)
    21 = @(7,8)
    22 = @(5,6)
    23 = @(8,8)
    24 = @(7,8)
    25 = @(8,8)
    26 = @(6,6)
    27 = @(7,8)
    28 = @(5,6)
    29 = @(8,8)
    30 = @(5,8)
    31 = @(8,9)
    32 = @(8,8)
    33 = @(8,8)
    34 = @(5,6)
    35 = @(4,7)
    36 = @(6,7)
    37 = @(7,8)
    38 = @(8,9)
    39 = @(8,9)
    40 = @(7,8)
    41 = @(7,8)
    42 = @(8,9)
    43 = @(6,7)
    44 = @(6,6)
    45 = @(6,6)
    46 = @(6,6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
